$d = $word.ActiveDocument

# Position a collapsed range at the very end of the document (right before
# the final section break), matching the spot where the new paragraphs
# belong per the diff.
$endRange = $d.Content
$endRange.Collapse(0)

# Build the OOXML for the two new paragraphs that get appended after the
# "... 15.50 Meer Art en scenes" paragraph:
#   1) a blank paragraph (nl-NL language mark only, no runs)
#   2) a paragraph with three separate runs describing the coroutine removal
$newParagraphsXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:rPr>
<w:lang w:val="nl-NL"/>
</w:rPr>
</w:pPr>
</w:p>
<w:p>
<w:pPr>
<w:rPr>
<w:lang w:val="nl-NL"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:lang w:val="nl-NL"/>
</w:rPr>
<w:t>14.00 start</w:t>
</w:r>
<w:r>
<w:rPr>
<w:lang w:val="nl-NL"/>
</w:rPr>
<w:t xml:space="preserve"> 15.00 Removed coroutines, werkt nu iets beter</w:t>
</w:r>
<w:r>
<w:rPr>
<w:lang w:val="nl-NL"/>
</w:rPr>
<w:t>. Nogsteeds de triggers werken niet</w:t>
</w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$endRange.InsertXML($newParagraphsXml)
